# "nueva tabla horas extras"
# Rename the "Horas Extra" sheet's hour-type columns (D1:I1) to be
# namespaced under "horas_extras/..." instead of the old bare names, and
# make "Horas Extra" the active/selected sheet (it was "Guía" before).

$wb = $excel.ActiveWorkbook

$wsHoras = $wb.Worksheets.Item("Horas Extra")

$wsHoras.Range("D1").Value = "horas_extras/horas_diurnas"
$wsHoras.Range("E1").Value = "horas_extras/horas_nocturnas"
$wsHoras.Range("F1").Value = "horas_extras/horas_diurnas_descanso"
$wsHoras.Range("G1").Value = "horas_extras/horas_nocturnas_descanso"
$wsHoras.Range("H1").Value = "horas_extras/horas_diurnas_asueto"
$wsHoras.Range("I1").Value = "horas_extras/horas_nocturnas_asueto"

# Make "Horas Extra" the active sheet/tab, with K5 selected (matches the
# author's working selection when they added the new columns).
$wsHoras.Activate()
$wsHoras.Range("K5").Select()
